$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 (shifts existing rows 22-34 down to 23-35),
# mirroring the weekly data refresh: a new reading (27-Aug-2021) is inserted
# at the top of the date-descending list, so every older reading slides down
# one row.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the latest weekly reading.
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value = 44435
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = 100112012
$ws.Cells.Item(22, 7).Value = "Espinaca"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 270
$ws.Cells.Item(22, 11).Value = 1800
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = 1900
$ws.Cells.Item(22, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 633
$ws.Cells.Item(22, 17).Value = 3
$ws.Cells.Item(22, 18).Value = "Hortaliza"

# Append one more weekly reading (24-Aug-2021) as a brand-new row at the end
# of the table (row 36, since the table now spans through row 35).
$ws.Cells.Item(36, 1).Value = 1
$ws.Cells.Item(36, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value = 44432
$ws.Cells.Item(36, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(36, 5).Value = 15
$ws.Cells.Item(36, 6).Value = 100112012
$ws.Cells.Item(36, 7).Value = "Espinaca"
$ws.Cells.Item(36, 8).Value = "Sin especificar"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 270
$ws.Cells.Item(36, 11).Value = 1800
$ws.Cells.Item(36, 12).Value = 2000
$ws.Cells.Item(36, 13).Value = 1900
$ws.Cells.Item(36, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(36, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value = 633
$ws.Cells.Item(36, 17).Value = 3
$ws.Cells.Item(36, 18).Value = "Hortaliza"

Write-Output ("Final dimension: " + $ws.UsedRange.Address())
